# Auto-generated edit script: updates profit/price columns (H-N) across 39 leve rows
# spanning all 8 sheets, per the commit diff against Sheets/Sagittarius_Profits.xlsx.
$wb = $excel.ActiveWorkbook

# --- ALC!row4 (Leve Item ID G4=5470) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 4200
$ws.Cells.Item(4, 9).Value = 3600
$ws.Cells.Item(4, 10).Value = 6000
$ws.Cells.Item(4, 11).Value = 3600
$ws.Cells.Item(4, 12).Value = 6000
$ws.Cells.Item(4, 13).Value = -3486
$ws.Cells.Item(4, 14).Value = -6228

# --- ALC!row9 (Leve Item ID G9=5487) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 3108.0667
$ws.Cells.Item(9, 9).Value = 4508.8887
$ws.Cells.Item(9, 11).Value = 4508.8887
$ws.Cells.Item(9, 13).Value = -4339.8887

# --- ALC!row43 (Leve Item ID G43=5472) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 4331.6665
$ws.Cells.Item(43, 10).Value = 4497.5
$ws.Cells.Item(43, 12).Value = 4497.5
$ws.Cells.Item(43, 14).Value = -4635.5

# --- ALC!row80 (Leve Item ID G80=12605) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(80, 8).Value = 58015.42
$ws.Cells.Item(80, 9).Value = 2619.8572
$ws.Cells.Item(80, 10).Value = 90329.5
$ws.Cells.Item(80, 11).Value = 7859.571599999999
$ws.Cells.Item(80, 12).Value = 270988.5
$ws.Cells.Item(80, 13).Value = -6861.571599999999
$ws.Cells.Item(80, 14).Value = -272984.5

# --- ALC!row83 (Leve Item ID G83=12605) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(83, 8).Value = 58015.42
$ws.Cells.Item(83, 9).Value = 2619.8572
$ws.Cells.Item(83, 10).Value = 90329.5
$ws.Cells.Item(83, 11).Value = 23578.7148
$ws.Cells.Item(83, 12).Value = 812965.5
$ws.Cells.Item(83, 13).Value = -18586.7148
$ws.Cells.Item(83, 14).Value = -822949.5

# --- ALC!row100 (Leve Item ID G100=19906) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(100, 8).Value = 3300.75
$ws.Cells.Item(100, 9).Value = 2283.4
$ws.Cells.Item(100, 10).Value = 4996.3335
$ws.Cells.Item(100, 11).Value = 2283.4
$ws.Cells.Item(100, 12).Value = 4996.3335
$ws.Cells.Item(100, 13).Value = -1742.4
$ws.Cells.Item(100, 14).Value = -6078.3335

# --- ALC!row116 (Leve Item ID G116=27778) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 4993.5
$ws.Cells.Item(116, 9).Value = 4993.5
$ws.Cells.Item(116, 11).Value = 4993.5
$ws.Cells.Item(116, 13).Value = -1551.5

# --- ALC!row137 (Leve Item ID G137=44013) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 1906.5
$ws.Cells.Item(137, 9).Value = 1607.2142
$ws.Cells.Item(137, 11).Value = 4821.642599999999
$ws.Cells.Item(137, 13).Value = -2271.642599999999

# --- ALC!row138 (Leve Item ID G138=44169) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 3813.1516
$ws.Cells.Item(138, 10).Value = 3896.2278
$ws.Cells.Item(138, 12).Value = 11688.6834
$ws.Cells.Item(138, 14).Value = -21968.6834

# --- ARM!row61 (Leve Item ID G61=43999) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 7025.6313
$ws.Cells.Item(61, 9).Value = 6222.778
$ws.Cells.Item(61, 10).Value = 7748.2
$ws.Cells.Item(61, 11).Value = 6222.778
$ws.Cells.Item(61, 12).Value = 7748.2
$ws.Cells.Item(61, 13).Value = -6010.778
$ws.Cells.Item(61, 14).Value = -8172.2

# --- ARM!row63 (Leve Item ID G63=12528) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 2638.5
$ws.Cells.Item(63, 9).Value = 2499.6
$ws.Cells.Item(63, 11).Value = 2499.6
$ws.Cells.Item(63, 13).Value = -1813.6

# --- ARM!row66 (Leve Item ID G66=12528) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(66, 8).Value = 2638.5
$ws.Cells.Item(66, 9).Value = 2499.6
$ws.Cells.Item(66, 11).Value = 12498
$ws.Cells.Item(66, 13).Value = -9066

# --- ARM!row76 (Leve Item ID G76=10679) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(76, 8).Value = 10000
$ws.Cells.Item(76, 10).Value = 10000
$ws.Cells.Item(76, 12).Value = 10000
$ws.Cells.Item(76, 14).Value = -10676

# --- ARM!row79 (Leve Item ID G79=10679) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(79, 8).Value = 10000
$ws.Cells.Item(79, 10).Value = 10000
$ws.Cells.Item(79, 12).Value = 10000
$ws.Cells.Item(79, 14).Value = -12340

# --- ARM!row136 (Leve Item ID G136=43999) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 7025.6313
$ws.Cells.Item(136, 9).Value = 6222.778
$ws.Cells.Item(136, 10).Value = 7748.2
$ws.Cells.Item(136, 11).Value = 18668.334
$ws.Cells.Item(136, 12).Value = 23244.6
$ws.Cells.Item(136, 13).Value = -16118.334
$ws.Cells.Item(136, 14).Value = -28344.6

# --- BSM!row99 (Leve Item ID G99=19943) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 3847.2856
$ws.Cells.Item(99, 9).Value = 2984
$ws.Cells.Item(99, 11).Value = 2984
$ws.Cells.Item(99, 13).Value = -1486

# --- BSM!row107 (Leve Item ID G107=27706) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 4999
$ws.Cells.Item(107, 9).Value = 4999
$ws.Cells.Item(107, 10).Value = 0
$ws.Cells.Item(107, 11).Value = 4999
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 13).Value = -3079
$ws.Cells.Item(107, 14).ClearContents()

# --- CRP!row31 (Leve Item ID G31=44023) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 11827.143
$ws.Cells.Item(31, 9).Value = 6265.6665
$ws.Cells.Item(31, 10).Value = 15998.25
$ws.Cells.Item(31, 11).Value = 6265.6665
$ws.Cells.Item(31, 12).Value = 15998.25
$ws.Cells.Item(31, 13).Value = -5970.6665
$ws.Cells.Item(31, 14).Value = -16588.25

# --- CRP!row34 (Leve Item ID G34=44023) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 11827.143
$ws.Cells.Item(34, 9).Value = 6265.6665
$ws.Cells.Item(34, 10).Value = 15998.25
$ws.Cells.Item(34, 11).Value = 6265.6665
$ws.Cells.Item(34, 12).Value = 15998.25
$ws.Cells.Item(34, 13).Value = -6063.6665
$ws.Cells.Item(34, 14).Value = -16402.25

# --- CRP!row99 (Leve Item ID G99=36198) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 27793.334
$ws.Cells.Item(99, 9).Value = 30947.111
$ws.Cells.Item(99, 11).Value = 30947.111
$ws.Cells.Item(99, 13).Value = -29449.111

# --- CRP!row100 (Leve Item ID G100=34388) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(100, 8).Value = 99997.5
$ws.Cells.Item(100, 10).Value = 99997.5
$ws.Cells.Item(100, 12).Value = 99997.5
$ws.Cells.Item(100, 14).Value = -102161.5

# --- CRP!row105 (Leve Item ID G105=19928) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(105, 8).Value = 3410.5264
$ws.Cells.Item(105, 9).Value = 2737.5
$ws.Cells.Item(105, 11).Value = 2737.5
$ws.Cells.Item(105, 13).Value = -990.5

# --- CRP!row107 (Leve Item ID G107=27689) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 0
$ws.Cells.Item(107, 9).Value = 0
$ws.Cells.Item(107, 10).Value = 0
$ws.Cells.Item(107, 11).Value = 0
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 13).ClearContents()
$ws.Cells.Item(107, 14).ClearContents()

# --- CRP!row126 (Leve Item ID G126=36198) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 27793.334
$ws.Cells.Item(126, 9).Value = 30947.111
$ws.Cells.Item(126, 11).Value = 92841.333
$ws.Cells.Item(126, 13).Value = -90371.333

# --- CRP!row141 (Leve Item ID G141=43345) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(141, 8).Value = 309412.8
$ws.Cells.Item(141, 10).Value = 309412.8
$ws.Cells.Item(141, 12).Value = 309412.8
$ws.Cells.Item(141, 14).Value = -319772.8

# --- CUL!row139 (Leve Item ID G139=44102) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(139, 8).Value = 3723
$ws.Cells.Item(139, 9).Value = 3563.375
$ws.Cells.Item(139, 10).Value = 5000
$ws.Cells.Item(139, 11).Value = 10690.125
$ws.Cells.Item(139, 12).Value = 15000
$ws.Cells.Item(139, 13).Value = -5550.125
$ws.Cells.Item(139, 14).Value = -25280

# --- GSM!row80 (Leve Item ID G80=12521) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 9397
$ws.Cells.Item(80, 9).Value = 3831.6667
$ws.Cells.Item(80, 10).Value = 11965.615
$ws.Cells.Item(80, 11).Value = 3831.6667
$ws.Cells.Item(80, 12).Value = 11965.615
$ws.Cells.Item(80, 13).Value = -2833.6667
$ws.Cells.Item(80, 14).Value = -13961.615

# --- GSM!row83 (Leve Item ID G83=12521) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(83, 8).Value = 9397
$ws.Cells.Item(83, 9).Value = 3831.6667
$ws.Cells.Item(83, 10).Value = 11965.615
$ws.Cells.Item(83, 11).Value = 19158.3335
$ws.Cells.Item(83, 12).Value = 59828.075
$ws.Cells.Item(83, 13).Value = -14166.3335
$ws.Cells.Item(83, 14).Value = -69812.075

# --- GSM!row96 (Leve Item ID G96=18261) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(96, 8).Value = 0
$ws.Cells.Item(96, 10).Value = 0
$ws.Cells.Item(96, 12).Value = 0
$ws.Cells.Item(96, 14).ClearContents()

# --- GSM!row102 (Leve Item ID G102=36169) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 0
$ws.Cells.Item(102, 9).Value = 0
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 11).Value = 0
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 13).ClearContents()
$ws.Cells.Item(102, 14).ClearContents()

# --- GSM!row126 (Leve Item ID G126=36184) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 2993.7778
$ws.Cells.Item(126, 9).Value = 2868
$ws.Cells.Item(126, 10).Value = 4000
$ws.Cells.Item(126, 11).Value = 8604
$ws.Cells.Item(126, 12).Value = 12000
$ws.Cells.Item(126, 13).Value = -6134
$ws.Cells.Item(126, 14).Value = -16940

# --- LTW!row22 (Leve Item ID G22=5277) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2229.6316
$ws.Cells.Item(22, 10).Value = 2500.4
$ws.Cells.Item(22, 12).Value = 2500.4
$ws.Cells.Item(22, 14).Value = -3090.4

# --- LTW!row27 (Leve Item ID G27=5277) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27, 8).Value = 2229.6316
$ws.Cells.Item(27, 10).Value = 2500.4
$ws.Cells.Item(27, 12).Value = 2500.4
$ws.Cells.Item(27, 14).Value = -2714.4

# --- LTW!row93 (Leve Item ID G93=19993) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 490.9091
$ws.Cells.Item(93, 9).Value = 433.14285
$ws.Cells.Item(93, 11).Value = 433.14285
$ws.Cells.Item(93, 13).Value = 814.85715

# --- LTW!row132 (Leve Item ID G132=44058) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 4537.4165
$ws.Cells.Item(132, 9).Value = 3994.4443
$ws.Cells.Item(132, 11).Value = 11983.3329
$ws.Cells.Item(132, 13).Value = -9453.332900000001

# --- WVR!row92 (Leve Item ID G92=18088) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(92, 8).Value = 51998.75
$ws.Cells.Item(92, 10).Value = 51998.75
$ws.Cells.Item(92, 12).Value = 51998.75
$ws.Cells.Item(92, 14).Value = -56990.75

# --- WVR!row96 (Leve Item ID G96=19977) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 9581.615
$ws.Cells.Item(96, 9).Value = 9429.333000000001
$ws.Cells.Item(96, 11).Value = 9429.333000000001
$ws.Cells.Item(96, 13).Value = -8056.333000000001

# --- WVR!row107 (Leve Item ID G107=27746) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 581.1429000000001
$ws.Cells.Item(107, 9).Value = 598.1429000000001
$ws.Cells.Item(107, 10).Value = 564.1429000000001
$ws.Cells.Item(107, 11).Value = 1794.4287
$ws.Cells.Item(107, 12).Value = 1692.4287
$ws.Cells.Item(107, 13).Value = 125.5712999999998
$ws.Cells.Item(107, 14).Value = -5532.4287

# --- WVR!row132 (Leve Item ID G132=44029) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 6278.75
$ws.Cells.Item(132, 9).Value = 6278.75
$ws.Cells.Item(132, 11).Value = 18836.25
$ws.Cells.Item(132, 13).Value = -16306.25
